# Add a new "2023" column (T) to the maternal-mortality-rate sheet,
# mirroring the existing per-region data columns (D..S = years 2007..2022).
#
# This reproduces the upstream commit that appended the 2023 figures:
#   - header cell T4 = 2023 (with the same bottom border as the rest of
#     the year-header row)
#   - T3 gets the same (empty, bottom-bordered) formatting as the rest of
#     the header-divider row
#   - T5..T14 get the 2023 values for each oblast / "Kyrgyz Republic" row,
#     formatted like the existing numeric columns (0.0)
#   - T9 ("Issyk-Kul oblast") has no 2023 figure, so it reuses the existing
#     "-" placeholder text, same as the other "missing data" cells in the
#     sheet (e.g. N9, M11, S14, ...)
#   - T14 (bottom row of the table) also keeps the thicker bottom border
#     used to close off the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BottomBorder($rng) {
    $b = $rng.Borders.Item(9)   # xlEdgeBottom
    $b.LineStyle = 1            # xlContinuous
    $b.Weight = -4138           # xlMedium
    $b.ColorIndex = 1           # automatic/black
}

# --- Year header (row 4) ---------------------------------------------------
$ws.Range("T4").Value = 2023

# --- Data rows (5-14): 2023 values per region ------------------------------
$ws.Range("T5").Value = 22.606300992622124
$ws.Range("T6").Value = 13.621194578764559
$ws.Range("T7").Value = 38.913029379337182
$ws.Range("T8").Value = 19.215987701767872
$ws.Range("T9").Value = "-"
$ws.Range("T10").Value = 19.9288256227758
$ws.Range("T11").Value = 48.820179007323027
$ws.Range("T12").Value = 18.458698661744346
$ws.Range("T13").Value = 7.704160246533128
$ws.Range("T14").Value = 29.197080291970806

# Match the "0.0" number formatting used by the rest of the data columns.
$ws.Range("T5:T14").NumberFormat = "0.0"

# Row 3 is the thin divider row above the year headers; row 4's header row
# and row 14 (the last data row) both close with a bottom border, same as
# the rest of the table - extend that border into the new column.
Set-BottomBorder($ws.Range("T3"))
Set-BottomBorder($ws.Range("T4"))
Set-BottomBorder($ws.Range("T14"))

# Row 1 header cell grew slightly taller to fit the extra column heading.
$ws.Rows.Item(1).RowHeight = 39.75

# Row 15 only had a lone formatted-but-empty cell (P15) before; it simply
# now spans one more column once the T-column cells above exist - no
# value needs to be written there.

# Move the selection back to the top of the sheet (matches the original
# workbook's default/active selection on this sheet).
$ws.Range("A1").Select()
